# Add example store/recall goals messages to the User_Initiated_Messages and
# Follow_Up_Messages sheets (pair programmed: Henry, Jay).

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("User_Initiated_Messages")
$ws4 = $wb.Worksheets.Item("Follow_Up_Messages")

# ---------------------------------------------------------------------------
# Sheet 3: User_Initiated_Messages
# Insert two new columns (Store, Recall) before the existing "Follow Ups"
# column, then add a new "recall goals" row.
# ---------------------------------------------------------------------------

$ws3.Columns("F:G").Insert()

$ws3.Range("F1").Value = "Store"
$ws3.Range("G1").Value = "Recall"

$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "recall goals"
$ws3.Range("C5").Value = "Message;Recall"
$ws3.Range("D5").Value = "Your goal was: [Goals]"
$ws3.Range("G5").Value = "Goals"

# ---------------------------------------------------------------------------
# Sheet 4: Follow_Up_Messages
# Same column insert, then mark the "benefits" and "limit days" prompts as
# Store messages (capturing Benefits / Goals), and replace the old "echo"
# placeholder text with a real follow-up message.
# ---------------------------------------------------------------------------

$ws4.Columns("F:G").Insert()

$ws4.Range("F1").Value = "Store"
$ws4.Range("G1").Value = "Recall"

# Row 2 ("Great, there's lots of benefits...") becomes a Store message that
# stores the answer under "Benefits".
$ws4.Range("C2").Value = "Message;Store"
$ws4.Range("F2").WrapText = $true
$ws4.Range("F2").Value = "Benefits"

# Row 5 ("Great, here are some quick tips...") gains the matching wrap-text
# formatting on the new blank columns, but is not itself a Store row.
$ws4.Range("F5").WrapText = $true
$ws4.Range("G5").WrapText = $true

# Row 7 ("Try to limit how many days you go gambling...") becomes a Store
# message that stores the answer under "Goals".
$ws4.Range("C7").Value = "Message;Store"
$ws4.Range("F7").WrapText = $true
$ws4.Range("F7").Value = "Goals"

# Row 9: replace the old "echo" placeholder with the real follow-up text,
# and give it the same blank wrap-text columns as the other note rows.
$ws4.Range("D9").Value = "I'll hold you to that goal!"
$ws4.Range("F9").WrapText = $true
$ws4.Range("G9").WrapText = $true

$wb.Save()
